$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the reviewdb data in row 26 -----------------------------------
# C26: sofershani9@gmail.com      -> rozend80@gmail.com
# D26: rotemzinger3@gmail.com     -> emmakrigel63@gmail.com  (trailing space kept)
$ws.Range("C26").Value = "rozend80@gmail.com"
$ws.Range("D26").Value = "emmakrigel63@gmail.com "

# --- Remove the stale hyperlink that used to sit on D26 --------------------
# The runtime only supports clearing *all* hyperlinks on a sheet in one shot
# (per-item Hyperlink.Delete() is a no-op here), so the only reliable way to
# drop a single hyperlink is to clear the whole collection and re-create
# every hyperlink except the one that must disappear (the one on D26, which
# pointed at the old "rotemzinger3@gmail.com" text).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C7"),  "mailto:jorjkluni03@gmail.com",      "", "", "jorjkluni03@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"),  "mailto:vikicrestina@gmail.com",     "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D8"),  "mailto:jorjkluni03@gmail.com",      "", "", "jorjkluni03@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C9"),  "mailto:nachumshayil@gmail.com",     "", "", "nachumshayil@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D9"),  "mailto:nachushay@gmail.com",        "", "", "nachushay@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:nevilgreen12@gmail.com",     "", "", "nevilgreen12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:vikicrestina@gmail.com",     "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:snizzvered@gmail.com",       "", "", "snizzvered@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:krigelron@gmail.com",        "", "", "krigelron@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:redvelvetmichael@gmail.com", "", "", "redvelvetmichael@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:veredsnir12@gmail.com",      "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:veredsnir12@gmail.com",      "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:kevinkors122@gmail.com",     "", "", "kevinkors122@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:freelancernachus@gmail.com", "", "", "freelancernachus@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:sm6502345@gmail.com",        "", "", "sm6502345@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:cybworking@gmail.com",       "", "", "cybworking@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:rontiddler560@gmail.com",    "", "", "rontiddler560@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D16"), "mailto:halachme@gmail.com",         "", "", "halachme@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D18"), "mailto:itaisenior@gmail.com",       "", "", "itaisenior@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C20"), "mailto:rotemzinger3@gmail.com",     "", "", "rotemzinger3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C21"), "mailto:sinuspai@gmail.com",         "", "", "sinuspai@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D22"), "mailto:rotemzinger3@gmail.com",     "", "", "rotemzinger3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D25"), "mailto:itaisenior@gmail.com",       "", "", "itaisenior@gmail.com")

# --- Update the saved view state -------------------------------------------
# Active cell / selection moves to C26:D26 (scroll position topLeftCell=A13
# is not something this host exposes for writing, so only the selection is
# reproduced here).
$ws.Range("C26:D26").Select()
